$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REN_5DH")

# Update the contact pair labels (order matters for shared string table layout)
$ws.Range("A2").Value = "GND (P1) to Cover Center (P2)"
$ws.Range("A4").Value = "GND (P1) to Cover Center (P2) (RE)"
$ws.Range("A3").Value = "GND (P1) to Cover Edge (P2)"
$ws.Range("A5").Value = "GND (P1) to Cover Edge (P2) (RE)"

# Make REN_5DH the active tab and select cell F8, as in the author's edit
$ws.Activate()
$ws.Range("F8").Select()
